$wb = $excel.ActiveWorkbook

# Source sheet whose error-code table is being reused for the new "CAL" sheet
$pt = $wb.Worksheets.Item("PT")

# Add the new sheet after the last existing sheet (AUTO_OBA) and name it "CAL"
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$cal = $wb.Worksheets.Add($null, $lastSheet)
$cal.Name = "CAL"

# Copy the CODE/Description table (rows 1-15) from PT into the new CAL sheet,
# preserving values, shared-string usage and styles
$pt.Range("A1:B15").Copy($cal.Range("A1:B15"))

# Match column B's width on the new sheet
$cal.Columns.Item(2).ColumnWidth = 30.33

# Update selections: PT no longer keeps its old A15 cell selection, CAL becomes
# the active sheet with A1:B15 selected
$pt.Range("A1:B15").Select()
$cal.Range("A1:B15").Select()
$cal.Activate()
